# Applies the "updated 4.0 files and mdl" revision to
# Fuel Prod Imp Exp Balancing Priorities.xlsx
#
#  - About sheet:  refresh the "last updated" date (C1) and scroll the
#                   view down a bit (so row 6 sits at the top).
#  - FPIEBP sheet: re-prioritize "hard coal" (row 3) production/imports/
#                   exports from 3/2/1 to 1/3/2, and move the active
#                   selection to E3.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsFpiebp = $wb.Worksheets.Item("FPIEBP")

# --- About sheet -----------------------------------------------------
# Update the date stamp in C1 (serial 45294 -> 45379, i.e. 2024-01-03 -> 2024-03-28)
$wsAbout.Range("C1").Value = 45379

# Scroll the view so row 6 becomes the top visible row (topLeftCell = A6)
$wsAbout.Activate()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1

# --- FPIEBP sheet ------------------------------------------------------
# "hard coal" row: production/imports/exports priorities change 3,2,1 -> 1,3,2
$wsFpiebp.Activate()
$wsFpiebp.Range("B3").Value = 1
$wsFpiebp.Range("C3").Value = 3
$wsFpiebp.Range("D3").Value = 2

# FPIEBP stays the selected/active tab, with the cursor moved to E3.
$wsFpiebp.Range("E3").Select() | Out-Null
